# revise mycompress for himeno
#
# Adds three new sz-mod error-bound columns (F:H) to the "himeno" summary
# table at row 15, fills in a new data row (17) for those columns (plus the
# existing B/D/E columns), repositions the three charts on the sheet, and
# updates the active sheet/selection so "himeno" (not "ping-pong") is the
# tab shown when the workbook is reopened.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("himeno")

# --- New header cells for row 15 (sz-mod-0.0001 / 0.001 / 0.01) ---------
$ws.Range("F15").Value = "sz-mod-0.0001"
$ws.Range("G15").Value = "sz-mod-0.001"
$ws.Range("H15").Value = "sz-mod-0.01"

# G15/H15 pick up a new (non-hyperlink) black-text font style.
$ws.Range("G15:H15").Font.Color = 0

# --- New data row 17 ------------------------------------------------------
$ws.Range("B17").Value = 59.392279000000002
$ws.Range("D17").Value = 62.361249999999998
$ws.Range("E17").Value = 59.768925000000003
$ws.Range("F17").Value = 59.029927999999998
$ws.Range("G17").Value = 60.659789000000004
$ws.Range("H17").Value = 61.856563999999999

# --- Reposition the three charts on "himeno" -------------------------------
$charts = $ws.ChartObjects()

$c1 = $charts.Item(1)
$c1.Left = 355.43874938484254
$c1.Top = 374.5791338582677
$c1.Width = 364.9791141732283
$c1.Height = 216.0000787401575

$c2 = $charts.Item(2)
$c2.Left = 564.2772533218504
$c2.Top = 116.0
$c2.Width = 267.273720472441
$c2.Height = 230.5683464566929

$c3 = $charts.Item(3)
$c3.Left = 17.5
$c3.Top = 360.5
$c3.Width = 292.548828125
$c3.Height = 231.0

# --- Active sheet / selection ---------------------------------------------
# "himeno" becomes the selected tab (was "ping-pong"); the cursor rests on
# the newly-added H17 cell.
$ws.Activate()
[void]$ws.Range("H17").Select()
